$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to Text format so numeric-looking strings
# (e.g. "381.92") are not silently converted to floating point numbers.
$priceCells = $ws.Range("D2:D51")
$priceCells.NumberFormat = "@"

$ws.Range('D2').Value = '51.446.99'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '2.983.41'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '381.92'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').Value = '103.83'
$ws.Range('E6').Value = '  +2.46%  '
$ws.Range('D7').Value = '0.544'
$ws.Range('E7').Value = '  +0.70%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.591'
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = '36.60'
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('D12').Value = '0.0857'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').Value = '3.458.16'
$ws.Range('E13').Value = '  +1.44%  '
$ws.Range('D14').Value = '18.41'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('D15').Value = '7.78'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').Value = '2.985.49'
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('D17').Value = '11.18'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').Value = '0.995'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').Value = '51.482.23'
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('D20').Value = '3.10'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('D22').Value = '0.0₃0963'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').Value = '70.23'
$ws.Range('E23').Value = '  +2.01%  '
$ws.Range('D24').Value = '267.22'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').Value = '3.23'
$ws.Range('E25').Value = '  +2.45%  '
$ws.Range('D26').Value = '7.84'
$ws.Range('E26').Value = '  -4.48%  '
$ws.Range('D27').Value = '7.32'
$ws.Range('E27').Value = '  -4.36%  '
$ws.Range('D28').Value = '0.169'
$ws.Range('E28').Value = '  +2.92%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '26.04'
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').Value = '10.37'
$ws.Range('E32').Value = '  +3.46%  '
$ws.Range('D33').Value = '34.73'
$ws.Range('E33').Value = '  +3.57%  '
$ws.Range('D34').Value = '51.37'
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').Value = '2.07'
$ws.Range('E35').Value = '  +1.11%  '
$ws.Range('E36').Value = '  -0.23%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('E38').Value = '  +3.01%  '
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '0.117'
$ws.Range('E40').Value = '  +0.78%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '1.84'
$ws.Range('E41').Value = '  +2.34%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '2.56'
$ws.Range('E42').Value = '  +2.90%  '
$ws.Range('D43').Value = '126.18'
$ws.Range('E43').Value = '  +4.60%  '
$ws.Range('D44').Value = '3.79'
$ws.Range('E44').Value = '  +11.27%  '
$ws.Range('D45').Value = '21.42'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('D46').Value = '2.04'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '2.36'
$ws.Range('E47').Value = '  +1.73%  '
$ws.Range('B48').Value = 'TheGraph'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D48').Value = '0.270'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('D49').Value = '2.025.57'
$ws.Range('E49').Value = '  +1.37%  '
$ws.Range('D50').Value = '3.283.31'
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('D51').Value = '0.0331'
$ws.Range('E51').Value = '  +0.54%  '

# Restore the default (Normal) style on the price column so no stray
# cell-level Text-format style reference is left behind.
$priceCells.Style = "Normal"
